$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final roster (6 students), rows 4-9, columns A (Name) and B (Student Number)
# Row 4: ALICE / A0000002U -> AARON TAN / A0123456U
$ws.Range("A4").Value = "AARON TAN"
$ws.Range("B4").Value = "A0123456U"

# Row 5: BENNY / A0000001U -> BU WEN JIN / A0123456M (shifted up)
$ws.Range("A5").Value = "BU WEN JIN"
$ws.Range("B5").Value = "A0123456M"

# Row 6: BU WEN JIN / A0123456M -> LAU XIN YEE / A0765432U (shifted up)
$ws.Range("A6").Value = "LAU XIN YEE"
$ws.Range("B6").Value = "A0765432U"

# Row 7: LAU XIN YEE / A0765432U -> LIM CHUN YONG / A0388443R
$ws.Range("A7").Value = "LIM CHUN YONG"
$ws.Range("B7").Value = "A0388443R"

# Row 8: RYAN LIM / A0000000U -> LIM JIA RUI RYAN / A0587314L
$ws.Range("A8").Value = "LIM JIA RUI RYAN"
$ws.Range("B8").Value = "A0587314L"

# Row 9: RYAN TAN / A0000004U -> NEO RUI EN MAYBELLINE / A0139345U
$ws.Range("A9").Value = "NEO RUI EN MAYBELLINE"
$ws.Range("B9").Value = "A0139345U"
